$d = $word.ActiveDocument

# Find the paragraph containing "Difficulty: Medium-High" and remove the entire paragraph
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Difficulty: Medium-High*") {
        $p.Range.Delete()
        break
    }
}
